# fix: prevent hidden columns from being labeled upon detecting changes
#
# For the group of rows 26-59, the "Änderung" (change) marker that had been
# written into column L no longer applies (those rows turned out to only
# differ in columns that are hidden, so they must not be flagged as
# changed). In addition, the rows where a new field-group starts (i.e. the
# first row for a given column-B label) were missing the grey "group
# header" shading that every other group-start row in the sheet already
# has (see rows 2, 9, 14 and 18) - that shading is restored here as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row numbers (1-based worksheet rows) that start a new field group and
# therefore need the full grey "group header" styling (font/fill/border
# copied from an already-correct group-header row such as row 2).
$fullChangeRows = @(26, 30, 37, 41, 48, 55, 57)

# All rows whose "Änderung" flag/style in column L must be cleared.
$allRows = 26..59

# Template ranges that already carry the styling we want to reproduce.
$headerTemplate = $ws.Range("A2:V2")
$lTemplate = $ws.Range("L21")

foreach ($r in $fullChangeRows) {
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $headerTemplate.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

foreach ($r in $allRows) {
    $lCell = $ws.Range("L" + $r)
    $lTemplate.Copy()
    $lCell.PasteSpecial($xlPasteFormats)
    $lCell.ClearContents()
}

$excel.CutCopyMode = 0
